$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($rng, [string]$val)
    # Force a Text number format before assigning so numeric-looking
    # strings (e.g. "216.23") are not silently coerced to a Double,
    # then restore the default style so no stray formatting is left
    # behind on cells that were plain (unstyled) text before.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-CellText $ws.Range("D2") "27.047.15"
Set-CellText $ws.Range("E2") "  +2.34%  "
Set-CellText $ws.Range("D3") "1.672.63"
Set-CellText $ws.Range("E3") "  +3.14%  "
Set-CellText $ws.Range("E4") "  +0.09%  "
Set-CellText $ws.Range("D5") "216.23"
Set-CellText $ws.Range("E5") "  +1.48%  "
Set-CellText $ws.Range("E6") "  +1.91%  "
Set-CellText $ws.Range("E7") "  +0.03%  "
Set-CellText $ws.Range("E8") "  +2.17%  "
Set-CellText $ws.Range("E9") "  +1.41%  "
Set-CellText $ws.Range("D10") "20.18"
Set-CellText $ws.Range("E10") "  +4.90%  "
Set-CellText $ws.Range("D11") "0.0891"
Set-CellText $ws.Range("E11") "  +5.03%  "
Set-CellText $ws.Range("D12") "1.909.60"
Set-CellText $ws.Range("E12") "  +3.28%  "
Set-CellText $ws.Range("D13") "1.676.40"
Set-CellText $ws.Range("E13") "  +3.49%  "
Set-CellText $ws.Range("E14") "  +1.54%  "
Set-CellText $ws.Range("D15") "65.78"
Set-CellText $ws.Range("E15") "  +3.03%  "
Set-CellText $ws.Range("E16") "  +2.27%  "
Set-CellText $ws.Range("D17") "27.070.09"
Set-CellText $ws.Range("E17") "  +2.39%  "
Set-CellText $ws.Range("D18") "235.33"
Set-CellText $ws.Range("E18") "  -0.36%  "
Set-CellText $ws.Range("E19") "  +1.66%  "
Set-CellText $ws.Range("D20") "7.71"
Set-CellText $ws.Range("E20") "  -1.67%  "
Set-CellText $ws.Range("E21") "  -0.12%  "
Set-CellText $ws.Range("D22") "4.47"
Set-CellText $ws.Range("E22") "  +3.53%  "
Set-CellText $ws.Range("D23") "9.27"
Set-CellText $ws.Range("E23") "  +1.51%  "
Set-CellText $ws.Range("D24") "2.23"
Set-CellText $ws.Range("E24") "  +1.47%  "
Set-CellText $ws.Range("D25") "145.29"
Set-CellText $ws.Range("E25") "  -1.23%  "
Set-CellText $ws.Range("E26") "  +1.18%  "
Set-CellText $ws.Range("E27") "  +0.50%  "
Set-CellText $ws.Range("E28") "  +2.13%  "
Set-CellText $ws.Range("E29") "  -0.12%  "
Set-CellText $ws.Range("E30") "  +0.22%  "
Set-CellText $ws.Range("E31") "  +1.66%  "
Set-CellText $ws.Range("E32") "  +2.03%  "
Set-CellText $ws.Range("D33") "1.450.29"
Set-CellText $ws.Range("E33") "  -4.76%  "
Set-CellText $ws.Range("E34") "  +5.46%  "
Set-CellText $ws.Range("E35") "  +5.99%  "
Set-CellText $ws.Range("E36") "  -0.40%  "
Set-CellText $ws.Range("D37") "0.573"
Set-CellText $ws.Range("E37") "  +0.84%  "
Set-CellText $ws.Range("D38") "0.896"
Set-CellText $ws.Range("E38") "  +7.33%  "
Set-CellText $ws.Range("E39") "  +1.73%  "
Set-CellText $ws.Range("D40") "6.08"
Set-CellText $ws.Range("E40") "  +3.36%  "
Set-CellText $ws.Range("E42") "  +10.27%  "
Set-CellText $ws.Range("E43") "  +3.31%  "
Set-CellText $ws.Range("D44") "66.06"
Set-CellText $ws.Range("E44") "  +5.30%  "
Set-CellText $ws.Range("D45") "1.819.00"
Set-CellText $ws.Range("E45") "  +3.30%  "
Set-CellText $ws.Range("D46") "0.779"
Set-CellText $ws.Range("E46") "  +2.08%  "
Set-CellText $ws.Range("D47") "90.32"
Set-CellText $ws.Range("E47") "  -0.35%  "
Set-CellText $ws.Range("E48") "  +1.81%  "
Set-CellText $ws.Range("B49") "BabyDogeCoin"
Set-CellText $ws.Range("C49") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-CellText $ws.Range("D49") "0.0₆0104"
Set-CellText $ws.Range("E49") "  -0.75%  "
Set-CellText $ws.Range("B50") "Algorand"
Set-CellText $ws.Range("C50") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-CellText $ws.Range("D50") "0.101"
Set-CellText $ws.Range("E50") "  +4.19%  "
Set-CellText $ws.Range("B51") "Cronos"
Set-CellText $ws.Range("C51") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-CellText $ws.Range("D51") "0.0508"
Set-CellText $ws.Range("E51") "  +1.37%  "
